$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 3589.1217388848
$ws.Range("D2").Value = 676.122019222639
